$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Imported")

$ws.Range("A3").Value = "News Item"
$ws.Range("B3").Value = "/bar-news-item-title"
$ws.Range("C3").Value = "Bar News Item Title"
$ws.Range("D3").Value = "bar news item description"
$ws.Range("E3").Value = "<p>bar <em>news item</em> text</p>"

$ws.Range("E3").Select()
